# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Country label swaps (the shared-strings table order changed, so the
#    rows that reference those shared strings now show different text).
# ---------------------------------------------------------------------------
$ws.Cells.Item(15,1).Value = "Austria"
$ws.Cells.Item(16,1).Value = "Belgica"

$ws.Cells.Item(81,1).Value = "Tunez"
$ws.Cells.Item(82,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(83,1).Value = "Albania"

# ---------------------------------------------------------------------------
# 2) Updated "last refreshed" timestamp banner.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 15:42"

# ---------------------------------------------------------------------------
# 3) Updated numeric figures (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes).
# ---------------------------------------------------------------------------

# Row 6 (originally Iran)
$ws.Cells.Item(6,2).Value = 68905
$ws.Cells.Item(6,3).Value = 694
$ws.Cells.Item(6,5).Value = 67440

# Row 11 (originally Corea del Sur)
$ws.Cells.Item(11,2).Value = 11575
$ws.Cells.Item(11,3).Value = 678
$ws.Cells.Item(11,5).Value = 11272
$ws.Cells.Item(11,7).Value = 19
$ws.Cells.Item(11,8).Value = 172

# Row 15 now shows Austria -> new figures for Austria
$ws.Cells.Item(15,2).Value = 6398
$ws.Cells.Item(15,3).Value = 810
$ws.Cells.Item(15,4).Value = 112
$ws.Cells.Item(15,5).Value = 6237
$ws.Cells.Item(15,6).Value = 28
$ws.Cells.Item(15,7).Value = 18
$ws.Cells.Item(15,8).Value = 49

# Row 16 now shows Belgica -> new figures for Belgica
$ws.Cells.Item(16,2).Value = 6235
$ws.Cells.Item(16,3).Value = 1298
$ws.Cells.Item(16,4).Value = 675
$ws.Cells.Item(16,5).Value = 5340
$ws.Cells.Item(16,6).Value = 605
$ws.Cells.Item(16,7).Value = 42
$ws.Cells.Item(16,8).Value = 220

# Row 23 (originally Malasia)
$ws.Cells.Item(23,2).Value = 2589
$ws.Cells.Item(23,3).Value = 35
$ws.Cells.Item(23,5).Value = 2522

# Row 33 (originally Indonesia)
$ws.Cells.Item(33,2).Value = 1130
$ws.Cells.Item(33,3).Value = 67
$ws.Cells.Item(33,5).Value = 1101

# Row 57
$ws.Cells.Item(57,5).Value = 435
$ws.Cells.Item(57,7).Value = 3
$ws.Cells.Item(57,8).Value = 7

# Row 65
$ws.Cells.Item(65,4).Value = 24
$ws.Cells.Item(65,5).Value = 257

# Row 67
$ws.Cells.Item(67,5).Value = 271
$ws.Cells.Item(67,7).Value = 1
$ws.Cells.Item(67,8).Value = 1

# Row 81 now shows Tunez -> new figures for Tunez
$ws.Cells.Item(81,2).Value = 200
$ws.Cells.Item(81,3).Value = 27
$ws.Cells.Item(81,5).Value = 192
$ws.Cells.Item(81,6).Value = 10
$ws.Cells.Item(81,7).Value = 1
$ws.Cells.Item(81,8).Value = 6

# Row 82 now shows Bosnia y Herzegovina -> new figures for Bosnia y Herzegovina
$ws.Cells.Item(82,2).Value = 185
$ws.Cells.Item(82,3).Value = 9
$ws.Cells.Item(82,4).Value = 2
$ws.Cells.Item(82,5).Value = 180
$ws.Cells.Item(82,6).Value = 1
$ws.Cells.Item(82,7).Value = 0
$ws.Cells.Item(82,8).Value = 3

# Row 83 now shows Albania -> new figures for Albania
$ws.Cells.Item(83,2).Value = 174
$ws.Cells.Item(83,3).Value = 28
$ws.Cells.Item(83,4).Value = 17
$ws.Cells.Item(83,5).Value = 151
$ws.Cells.Item(83,6).Value = 3

# Row 86
$ws.Cells.Item(86,4).Value = 20
$ws.Cells.Item(86,5).Value = 133
